$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: "_old" -> "_FV2310", "_new" -> "_FV2404"
for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cur = $cell.Value()
    $cell.Value = ($cur -replace '_old$', '_FV2310')
}
for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cur = $cell.Value()
    $cell.Value = ($cur -replace '_new$', '_FV2404')
}

# Add an Excel Table (ListObject) over A1:U57
$rng = $ws.Range("A1:U57")
$tbl = $ws.ListObjects.Add(1, $rng, [System.Type]::Missing, 1)
$tbl.Name = "Table1"

# Freeze header row (split pane above row 2)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
